# Generate Report for Handoff
# Updates the "Latest Handoff Date/Datetime" values for the
# 3e724b28-b7b6-496f-8e81-9d2ce671322e.md file row (row 6) across all
# three sheets, reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date"
$overview.Range("D6").Value = "2016-31-19 18:31:21"

# zh-cn sheet: column E = "Latest Handoff Datetime"
$zhcn.Range("E6").Value = "2016-03-19 18:31:18"

# de-de sheet: column E = "Latest Handoff Datetime"
$dede.Range("E6").Value = "2016-03-19 18:31:21"
